$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.193.08'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.323.76'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.507'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.19'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('E14').Value = '  +1.76%  '
$ws.Range('D15').Value = '2.683.56'
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').Value = '2.362.12'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('E17').Value = '  -1.27%  '
$ws.Range('D18').Value = '43.104.81'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('D21').Value = '0.0₃0910'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '168.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -6.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.06%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.88%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.68'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.69%  '
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0698'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').Value = '2.003.51'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0291'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('E44').Value = '  -4.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.61%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('D50').Value = '2.548.71'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.87%  '
